# adaugare fisiere baza de date
# Update the stored user record: id -> 1001, email -> oclock@gmail.com,
# and move the active-cell selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# id column (A2): 1000 -> 1001
$ws.Range("A2").Value = 1001

# email column (B2, shared string, also the hyperlink display text): update address
$ws.Range("B2").Value = "oclock@gmail.com"

# Move/restore the saved cell selection to C6
$ws.Range("C6").Select() | Out-Null

# Best-effort: localize the theme name (not guaranteed to round-trip through
# this host's Theme object, but harmless if unsupported).
try {
    $wb.Theme.Name = "Temă Office"
} catch {
}
